# Apply updated cryptocurrency price/volume data to the active worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '66.258.51'
$ws.Range('E2').Value = '  +3.13%  '
$ws.Range('D3').Value = '3.247.66'
$ws.Range('E3').Value = '  +6.89%  '
$ws.Range('E4').Value = '  +0.09%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '581.59'
$ws.Range('E5').Value = '  +4.72%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '152.14'
$ws.Range('E6').Value = '  +7.49%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.999'
$ws.Range('E7').Value = '  -0.15%  '
$ws.Range('D8').Value = '3.239.31'
$ws.Range('E8').Value = '  +7.01%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.515'
$ws.Range('E9').Value = '  +5.84%  '
$ws.Range('E10').Value = '  +9.88%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.164'
$ws.Range('E11').Value = '  +6.88%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.487'
$ws.Range('E12').Value = '  +5.08%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '37.84'
$ws.Range('E13').Value = '  +3.99%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.0000235'
$ws.Range('E14').Value = '  +6.61%  '
$ws.Range('D15').Value = '3.778.58'
$ws.Range('E15').Value = '  +7.22%  '
$ws.Range('D16').Value = '66.423.58'
$ws.Range('E16').Value = '  +3.29%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '553.29'
$ws.Range('E17').Value = '  +13.61%  '
$ws.Range('D18').Value = '3.258.22'
$ws.Range('E18').Value = '  +7.05%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.114'
$ws.Range('E19').Value = '  +2.86%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '7.13'
$ws.Range('E20').Value = '  +6.29%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '14.50'
$ws.Range('E21').Value = '  +6.24%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.744'
$ws.Range('E22').Value = '  +8.39%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '7.86'
$ws.Range('E23').Value = '  +10.66%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '13.57'
$ws.Range('E24').Value = '  +7.69%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '81.56'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.998'
$ws.Range('E26').Value = '  +0.02%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '9.24'
$ws.Range('E27').Value = '  +18.90%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.98'
$ws.Range('E28').Value = '  +8.94%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.25'
$ws.Range('E29').Value = '  +7.05%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '27.77'
$ws.Range('E30').Value = '  +7.24%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '2.77'
$ws.Range('E31').Value = '  +6.74%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.00'
$ws.Range('E32').Value = '  +0.20%  '
$ws.Range('E33').Value = '  +5.56%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '565.09'
$ws.Range('E34').Value = '  +8.09%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '5.68'
$ws.Range('E35').Value = '  +4.18%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '6.38'
$ws.Range('E36').Value = '  +7.17%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '55.34'
$ws.Range('E37').Value = '  +5.51%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.0453'
$ws.Range('E38').Value = '  +12.11%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.0863'
$ws.Range('E39').Value = '  +8.23%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.130'
$ws.Range('E40').Value = '  +7.86%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '2.99'
$ws.Range('E41').Value = '  +9.66%  '
$ws.Range('D42').Value = '3.215.03'
$ws.Range('E42').Value = '  +11.54%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '8.64'
$ws.Range('E43').Value = '  +3.82%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.282'
$ws.Range('E44').Value = '  +15.44%  '
$ws.Range('E45').Value = '  +10.61%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '26.51'
$ws.Range('E46').Value = '  +6.47%  '
$ws.Range('D48').Value = '0.0₃0557'
$ws.Range('E48').Value = '  +4.61%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '125.24'
$ws.Range('E49').Value = '  +3.81%  '
$ws.Range('E50').Value = '  +4.34%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '2.21'
$ws.Range('E51').Value = '  +8.69%  '
